# Lab-Machines.xlsx update ("Add files via upload"):
#  - Column A (User No placeholders) replaced with the real trainee names
#  - Two more trainees appended as new rows (13 and 14), name only
#  - A thin border is drawn around the whole table (A1:G14)
#  - A few columns are widened to fit the longer names / values
#  - Zoom level and active selection changed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Replace the "User N" placeholders in column A with real names ----
$names = @(
    "Bharath P",
    "Subrat Kumar Das",
    "Radhika T",
    "Rajeshwari A",
    "Subburaj A",
    "Rohit Vignesh",
    "Balamurugan G",
    "Srinivas K",
    "Ashish Gupta",
    "Suresh",
    "Saravanan Rajamanickam"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

# ---- 2. Append two more trainees as new rows, name only ----
$ws.Cells.Item(13, 1).Value = "Abdul Razack"
$ws.Cells.Item(14, 1).Value = "Vijay"

# ---- 3. Draw a thin border around the whole used range ----
$tableRange = $ws.Range("A1:G14")
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2

# ---- 4. Center the two new name cells like the rest of column A ----
$newNamesRange = $ws.Range("A13:A14")
$newNamesRange.VerticalAlignment = -4108
$newNamesRange.HorizontalAlignment = -4108

# ---- 5. Resize columns to fit the new (longer) names / widen a few ----
$ws.Columns.Item(1).ColumnWidth = 22.736979166666668
$ws.Columns.Item(2).ColumnWidth = 19.166666666666668
$ws.Columns.Item(3).ColumnWidth = 14.592447916666666
$ws.Columns.Item(5).ColumnWidth = 17.307291666666668
$ws.Columns.Item(6).ColumnWidth = 14.166666666666666

# ---- 6. Zoom + selection, like the saved workbook ----
$excel.ActiveWindow.Zoom = 130
$ws.Range("B21").Select()
